$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: prepare new row 132 by copying the cell formats (styles) of row 131 ---
# Column A (id) uses style index 1 (bold/border), column E (Date) uses style index 2 (date format).
# Copying formats only (not values) lets the engine reuse the existing style indices
# instead of creating brand-new duplicate style entries in styles.xml.
$ws.Cells.Item(131, 1).Copy() | Out-Null
$ws.Cells.Item(132, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(131, 5).Copy() | Out-Null
$ws.Cells.Item(132, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Step 2: push the existing match (row 131) data down into the new row 132 ---
# This is the same fixture (Mohun Bagan SG vs Mumbai City FC, match id 7749764) that used
# to be row 131 (id=129); it now becomes row 132 with id=130, and two of its closing-line
# odds (oddAHOver / oddAHUnder -> columns U/V) were refreshed.
$ws.Cells.Item(132, 1).Value2 = 130
$ws.Cells.Item(132, 2).Value2 = 7749764
$ws.Cells.Item(132, 3).Value2 = "India Super League"
$ws.Cells.Item(132, 4).Value2 = "India Super League"
$ws.Cells.Item(132, 5).Value2 = 45397.45833333334
$ws.Cells.Item(132, 6).Value2 = "Mohun Bagan SG"
$ws.Cells.Item(132, 7).Value2 = "Mumbai City FC"
$ws.Cells.Item(132, 11).Value2 = 2.5
$ws.Cells.Item(132, 12).Value2 = 3.4
$ws.Cells.Item(132, 13).Value2 = 2.7
$ws.Cells.Item(132, 14).Value2 = 2.5
$ws.Cells.Item(132, 15).Value2 = 3.4
$ws.Cells.Item(132, 16).Value2 = 2.7
$ws.Cells.Item(132, 17).Value2 = 0
$ws.Cells.Item(132, 18).Value2 = 1.825
$ws.Cells.Item(132, 19).Value2 = 1.975
$ws.Cells.Item(132, 20).Value2 = 2.75
$ws.Cells.Item(132, 21).Value2 = 1.95
$ws.Cells.Item(132, 22).Value2 = 1.85
$ws.Cells.Item(132, 23).Value2 = 0
$ws.Cells.Item(132, 24).Value2 = 0
$ws.Cells.Item(132, 25).Value2 = 0
$ws.Cells.Item(132, 26).Value2 = 0
$ws.Cells.Item(132, 27).Value2 = 0

# --- Step 3: overwrite row 131 with the newly-added fixture ---
# Northeast United 3-0 Odisha FC, played 2024-04-13 (id=129, match id 7749775).
$ws.Cells.Item(131, 1).Value2 = 129
$ws.Cells.Item(131, 2).Value2 = 7749775
$ws.Cells.Item(131, 3).Value2 = "India Super League"
$ws.Cells.Item(131, 4).Value2 = "India Super League"
$ws.Cells.Item(131, 5).Value2 = 45395.45833333334
$ws.Cells.Item(131, 6).Value2 = "Northeast United"
$ws.Cells.Item(131, 7).Value2 = "Odisha FC"
$ws.Cells.Item(131, 8).Value2 = 3
$ws.Cells.Item(131, 9).Value2 = 0
$ws.Cells.Item(131, 10).Value2 = "H"
$ws.Cells.Item(131, 11).Value2 = 2.8
$ws.Cells.Item(131, 12).Value2 = 3.5
$ws.Cells.Item(131, 13).Value2 = 2.2
$ws.Cells.Item(131, 14).Value2 = 2.6
$ws.Cells.Item(131, 15).Value2 = 3.5
$ws.Cells.Item(131, 16).Value2 = 2.375
$ws.Cells.Item(131, 17).Value2 = 0
$ws.Cells.Item(131, 18).Value2 = 1.975
$ws.Cells.Item(131, 19).Value2 = 1.825
$ws.Cells.Item(131, 20).Value2 = 3
$ws.Cells.Item(131, 21).Value2 = 2.025
$ws.Cells.Item(131, 22).Value2 = 1.775
$ws.Cells.Item(131, 23).Value2 = 1.6
$ws.Cells.Item(131, 24).Value2 = -1
$ws.Cells.Item(131, 25).Value2 = -1
$ws.Cells.Item(131, 26).Value2 = 0.9750000000000001
$ws.Cells.Item(131, 27).Value2 = -1
$ws.Cells.Item(131, 28).Value2 = 0
$ws.Cells.Item(131, 29).Value2 = -0

Write-Host "Rows 131-132 updated"
